# daily auto push: 2026-01-24 18:44 UTC
# Two new observation rows for 2026/01/24 (Sat) 23:00 and 2026/01/25 (Sun) 02:00
# are inserted right after the existing 2026/01/24 19:00 row (row 705),
# pushing every subsequent row (old 706..747) down by two (new 708..749).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 706:707 (everything currently at/after 706 shifts down by 2)
$ws.Range("A706:A707").EntireRow.Insert()

# New row 706: 2026/01/24, 土, 23, 201
$ws.Range("A706").NumberFormat = "@"
$ws.Range("A706").Value = "2026/01/24"
$ws.Range("A706").Style = "Normal"
$ws.Range("B706").Value = "土"
$ws.Range("C706").Value = 23
$ws.Range("D706").Value = 201

# New row 707: 2026/01/25, 日, 2, 201
$ws.Range("A707").NumberFormat = "@"
$ws.Range("A707").Value = "2026/01/25"
$ws.Range("A707").Style = "Normal"
$ws.Range("B707").Value = "日"
$ws.Range("C707").Value = 2
$ws.Range("D707").Value = 201
